$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the title row (row 1: "Hate Crimes - On-campus Student Housing Facilities").
# This shifts every row up by one, turning the old header row (row 2) into row 1
# and the old data rows (3..55) into rows 2..54.
$ws.Rows.Item(1).Delete()

# Rename the worksheet/tab.
$ws.Name = "hate-crimes-on-campus-student-h"

# Re-case / correct the header row text (now row 1 after the deletion above).
$ws.Range("A1").Value = "Survey Year"
$ws.Range("B1").Value = "UnitID"
$ws.Range("C1").Value = "Institution Name"
$ws.Range("D1").Value = "Campus ID"
$ws.Range("E1").Value = "Campus Name"
$ws.Range("F1").Value = "Institution Size"
$ws.Range("G1").Value = "Murder/Non-Negligent Manslaughter"
$ws.Range("H1").Value = "Murder/Non-Negligent Manslaughter - Race"
$ws.Range("I1").Value = "Murder/Non-Negligent Manslaughter - Religion"
$ws.Range("J1").Value = "Murder/Non-Negligent Manslaughter - Sexual Orientation"
$ws.Range("K1").Value = "Murder/Non-Negligent Manslaughter - Gender"
$ws.Range("L1").Value = "Murder/Non-Negligent Manslaughter - Disability"
$ws.Range("M1").Value = "Murder/Non-Negligent Manslaughter - Ethnicity/National Origin"
$ws.Range("N1").Value = "Negligent Manslaughter"
$ws.Range("O1").Value = "Negligent Manslaughter - Race"
$ws.Range("P1").Value = "Negligent Manslaughter - Religion"
$ws.Range("Q1").Value = "Negligent Manslaughter - Sexual Orientation"
$ws.Range("R1").Value = "Negligent Manslaughter - Gender"
$ws.Range("S1").Value = "Negligent Manslaughter - Disability"
$ws.Range("T1").Value = "Negligent Manslaughter - Ethnicity/National Origin"
$ws.Range("U1").Value = "Sex Offenses - Forcible"
$ws.Range("V1").Value = "Sex Offenses - Forcible - Race"
$ws.Range("W1").Value = "Sex Offenses - Forcible - Religion"
$ws.Range("X1").Value = "Sex Offenses - Forcible - Sexual Orientation"
$ws.Range("Y1").Value = "Sex Offenses - Forcible - Gender"
$ws.Range("Z1").Value = "Sex Offenses - Forcible - Disability"
$ws.Range("AA1").Value = "Sex Offenses - Forcible - Ethnicity/National Origin"
$ws.Range("AB1").Value = "Sex Offenses - Non-Forcible"
$ws.Range("AC1").Value = "Sex Offenses - Non-Forcible -Race"
$ws.Range("AD1").Value = "Sex Offenses - Non-Forcible - Religion"
$ws.Range("AE1").Value = "Sex Offenses - Non-Forcible - Sexual Orientation"
$ws.Range("AF1").Value = "Sex Offenses - Non-Forcible - Gender"
$ws.Range("AG1").Value = "Sex Offenses - Non-Forcible - Disability"
$ws.Range("AH1").Value = "Sex Offenses - Non-Forcible - Ethnicity/National Origin"
$ws.Range("AI1").Value = "Robbery"
$ws.Range("AJ1").Value = "Robbery - Race"
$ws.Range("AK1").Value = "Robbery - Religion"
$ws.Range("AL1").Value = "Robbery - Sexual Orientation"
$ws.Range("AM1").Value = "Robbery - Gender"
$ws.Range("AN1").Value = "Robbery - Disability"
$ws.Range("AO1").Value = "Robbery - Ethnicity/National Origin"
$ws.Range("AP1").Value = "Aggravated Assault"
$ws.Range("AQ1").Value = "Aggravated Assault - Race"
$ws.Range("AR1").Value = "Aggravated Assault - Religion"
$ws.Range("AS1").Value = "Aggravated Assault - Sexual Orientation"
$ws.Range("AT1").Value = "Aggravated Assault - Gender"
$ws.Range("AU1").Value = "Aggravated Assault - Disability"
$ws.Range("AV1").Value = "Aggravated Assault - Ethnicity/National Origin"
$ws.Range("AW1").Value = "Burglary"
$ws.Range("AX1").Value = "Burglary - Race"
$ws.Range("AY1").Value = "Burglary - Religion"
$ws.Range("AZ1").Value = "Burglary - Sexual Orientation"
$ws.Range("BA1").Value = "Burglary - Gender"
$ws.Range("BB1").Value = "Burglary - Disability"
$ws.Range("BC1").Value = "Burglary - Ethnicity/National Origin"
$ws.Range("BD1").Value = "Motor Vehicle Theft"
$ws.Range("BE1").Value = "Motor Vehicle Theft - Race"
$ws.Range("BF1").Value = "Motor Vehicle Theft - Religion"
$ws.Range("BG1").Value = "Motor Vehicle Theft - Sexual Orientation"
$ws.Range("BH1").Value = "Motor Vehicle Theft - Gender"
$ws.Range("BI1").Value = "Motor Vehicle Theft - Disability"
$ws.Range("BJ1").Value = "Motor Vehicle Theft - Ethnicity/National Origin"
$ws.Range("BK1").Value = "Arson"
$ws.Range("BL1").Value = "Arson - Race"
$ws.Range("BM1").Value = "Arson - Religion"
$ws.Range("BN1").Value = "Arson - Sexual Orientation"
$ws.Range("BO1").Value = "Arson - Gender"
$ws.Range("BP1").Value = "Arson - Disability"
$ws.Range("BQ1").Value = "Arson - Ethnicity/National Origin"
$ws.Range("BR1").Value = "Simple Assault"
$ws.Range("BS1").Value = "Simple Assault - Race"
$ws.Range("BT1").Value = "Simple Assault - Religion"
$ws.Range("BU1").Value = "Simple Assault - Sexual Orientation"
$ws.Range("BV1").Value = "Simple Assault - Gender"
$ws.Range("BW1").Value = "Simple Assault - Disability"
$ws.Range("BX1").Value = "Simple Assault - Ethnicity/National Origin"
$ws.Range("BY1").Value = "Larceny-Theft"
$ws.Range("BZ1").Value = "Larceny-Theft - Race"
$ws.Range("CA1").Value = "Larceny-Theft - Religion"
$ws.Range("CB1").Value = "Larceny-Theft - Sexual Orientation"
$ws.Range("CC1").Value = "Larceny-Theft - Gender"
$ws.Range("CD1").Value = "Larceny-Theft - Disability"
$ws.Range("CE1").Value = "Larceny-Theft - Ethnicity/National Origin"
$ws.Range("CF1").Value = "Intimidation"
$ws.Range("CG1").Value = "Intimidation - Race"
$ws.Range("CH1").Value = "Intimidation - Religion"
$ws.Range("CI1").Value = "Intimidation - Sexual Orientation"
$ws.Range("CJ1").Value = "Intimidation - Gender"
$ws.Range("CK1").Value = "Intimidation - Disability"
$ws.Range("CL1").Value = "Intimidation - Ethnicity/National Origin"
$ws.Range("CM1").Value = "Destruction/Damage/Vandalism of Property"
$ws.Range("CN1").Value = "Destruction/Damage/Vandalism of Property - Race"
$ws.Range("CO1").Value = "Destruction/Damage/Vandalism of Property - Religion"
$ws.Range("CP1").Value = "Destruction/Damage/Vandalism of Property - Sexual Orientation"
$ws.Range("CQ1").Value = "Destruction/Damage/Vandalism of Property - Gender"
$ws.Range("CR1").Value = "Destruction/Damage/Vandalism of Property - Disability"
$ws.Range("CS1").Value = "Destruction/Damage/Vandalism of Property - Ethnicity/National Origin"
